# "Modificacion en reportes y fechas"
#
# - Adds a new pending-task row about missing decimals in reports.
# - Marks a couple of existing tasks ("revisar los TODO" and the validation
#   task) as 100% done with a responsible person assigned.
# - Flips a few "en proceso" (in-progress) status cells over to "100%"
#   (percentage) values now that those tasks are finished.
# - Moves the active selection/viewport back to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: "revisar los TODO" -> Agustina, 100%
$ws.Range("B9").Value = "Agustina"
$ws.Range("C9").NumberFormat = "0%"
$ws.Range("C9").Value = 1

# Row 28: "Validacion en creacion de cuota..." -> Lucas, 100%
$ws.Range("B28").Value = "Lucas"
$ws.Range("C28").NumberFormat = "0%"
$ws.Range("C28").Value = 1

# Rows 30, 35, 36: "en proceso" text status -> 100% numeric status
$ws.Range("C30").NumberFormat = "0%"
$ws.Range("C30").Value = 1

$ws.Range("C35").NumberFormat = "0%"
$ws.Range("C35").Value = 1

$ws.Range("C36").NumberFormat = "0%"
$ws.Range("C36").Value = 1

# New row 41: new pending task about missing decimals in reports
$ws.Range("A41").Value = "en los reportes faltan los decimales"
$ws.Range("B41").Value = "Agustina"
$ws.Range("C41").NumberFormat = "0%"
$ws.Range("C41").Value = 1

# Reset view: drop the scrolled-down viewport and move the selection to C10
[void]$ws.Range("C10").Select()
